$wb = $excel.ActiveWorkbook

# 1. Create "Challenge 8" by copying "Challenge 7" (this preserves all formatting,
#    column widths, merged cells, styles, etc.) and placing it right after it.
$ws7 = $wb.Worksheets.Item("Challenge 7")
$ws7.Copy([System.Reflection.Missing]::Value, $ws7)
$ws8 = $wb.Worksheets.Item($ws7.Index + 1)
$ws8.Name = "Challenge 8"

# 2. Remove the six rows that belonged to the middle of the old "objectives" list
#    (rows 17-22). This leaves the old rows 23-26 shifted up to become the new
#    rows 17-20, which already carry the correct cell styles/borders for this
#    table layout (including the thick-bottom border on the last row).
$ws8.Rows("17:22").Delete()

# 3. Update the text content for Challenge 8.
$ws8.Range("A2").Value = "Challenge 8: Tightening database security"
$ws8.Range("A16").Value = "Applying security to the database"

$ws8.Range("B18").Value = "Successfully encrypted all columns in the [Accounts] table"
$ws8.Range("B17").Value = "Successfully encrypted the [Users] table"
$ws8.Range("B19").Value = "Successfully encrypted the [Tranactions].[AccountId] column"
$ws8.Range("B20").Value = "Application connectivity to the database is conducted via a managed principal"

# 4. The last row's objective text needs to wrap (it's noticeably longer), so turn
#    wrapping on for that cell while the existing border/fill stay intact.
$ws8.Range("B20").WrapText = $true

# 5. Update the points-possible values for this challenge.
$ws8.Range("D16").Value = 16
$ws8.Range("E16").Value = 13

# 6. Make "Challenge 8" the active/visible tab.
$ws8.Activate()
$ws8.Range("E26").Select()
